$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from similar "filled" rows onto the new row-10 cells
$ws.Range("E6").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)

# Fill in the new contact's data (row 10)
$ws.Range("A10").Value = 14000523
$ws.Range("B10").Value = 7265
$ws.Range("C10").Value = "علی محمدی"
$ws.Range("D10").Value = "ایلام"
$ws.Range("F10").Value = "08433303086"

# G10 is no longer blank-with-dropdown; drop it from the validated range
$ws.Range("G10").Validation.Delete()

# Move the active selection to H8
$ws.Range("H8").Select()
